# InsideBet Data: Automatizado
# The fixture that was the next upcoming match (Gameweek 23, Fri 2026-02-20,
# Brest vs Marseille) has passed / been removed from the "proximos partidos"
# feed. Delete its row entirely - Excel's native row-delete shifts every
# following row up by one, re-numbering week/date/time/home/away/venue data
# accordingly and dropping the now-unused final row (142).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("24:24").Delete()
